$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- colour + layout constants -------------------------------------------------
$RED    = 255        # RGB(255,0,0)   -> wrong
$BLUE   = 16711680    # RGB(0,0,255)   -> "correct answer" column
$GREEN  = 32768        # RGB(0,128,0)   -> right

function Format-Plain($cell) {
    # font2 Century/12, boxed, centered  (style "4" in the target workbook)
    $cell.Font.Name = "Century"
    $cell.Font.Size = 12
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
}

function Format-Colored($cell, $color) {
    # Century/12 in a given colour, boxed, centered (styles "5"/"6"/"7")
    $cell.Font.Name = "Century"
    $cell.Font.Size = 12
    $cell.Font.Color = $color
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
}

function Format-Bold($cell) {
    # font3 Century/12 bold, boxed, centered (style "8")
    $cell.Font.Name = "Century"
    $cell.Font.Size = 12
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
}

# ---- row 8 : blank spacer row ---------------------------------------------------
# ".Style = 'Normal'" forces the (otherwise empty) cells to materialise as real,
# default-styled records so the row exists - without fabricating a new style.
$ws.Range("A8:E8").Style = "Normal"

# ---- row 9 : column headers for the score box -----------------------------------
$ws.Cells.Item(9,2).Value = "Right"
$ws.Cells.Item(9,3).Value = "Wrong"
$ws.Cells.Item(9,4).Value = "Not Attempt"
$ws.Cells.Item(9,5).Value = "Max"
Format-Plain $ws.Cells.Item(9,1)
Format-Plain $ws.Cells.Item(9,2)
Format-Plain $ws.Cells.Item(9,3)
Format-Plain $ws.Cells.Item(9,4)
Format-Plain $ws.Cells.Item(9,5)

# ---- row 10 : No. -----------------------------------------------------------
$ws.Cells.Item(10,1).Value = "No."
$ws.Cells.Item(10,2).Value = 18
$ws.Cells.Item(10,3).Value = 7
$ws.Cells.Item(10,4).Value = 3
$ws.Cells.Item(10,5).Value = 28
Format-Plain   $ws.Cells.Item(10,1)
Format-Colored $ws.Cells.Item(10,2) $GREEN
Format-Colored $ws.Cells.Item(10,3) $RED
Format-Plain   $ws.Cells.Item(10,4)
Format-Plain   $ws.Cells.Item(10,5)

# ---- row 11 : Marking --------------------------------------------------------
$ws.Cells.Item(11,1).Value = "Marking"
$ws.Cells.Item(11,2).Value = 5
$ws.Cells.Item(11,3).Value = -1
$ws.Cells.Item(11,4).Value = 0
Format-Plain   $ws.Cells.Item(11,1)
Format-Colored $ws.Cells.Item(11,2) $GREEN
Format-Colored $ws.Cells.Item(11,3) $RED
Format-Plain   $ws.Cells.Item(11,4)
Format-Plain   $ws.Cells.Item(11,5)

# ---- row 12 : Total ----------------------------------------------------------
$ws.Cells.Item(12,1).Value = "Total"
$ws.Cells.Item(12,2).Value = 90
$ws.Cells.Item(12,3).Value = -7
$ws.Cells.Item(12,5).Value = "83/140"
Format-Plain   $ws.Cells.Item(12,1)
Format-Colored $ws.Cells.Item(12,2) $GREEN
Format-Colored $ws.Cells.Item(12,3) $RED
Format-Plain   $ws.Cells.Item(12,4)
Format-Colored $ws.Cells.Item(12,5) $BLUE

# ---- rows 13-14 : blank spacer rows ---------------------------------------------
$ws.Range("A13:E13").Style = "Normal"
$ws.Range("A14:E14").Style = "Normal"

# ---- row 15 : "Student Ans" / "Correct Ans" headers (two side-by-side tables) ---
$ws.Cells.Item(15,1).Value = "Student Ans"
$ws.Cells.Item(15,2).Value = "Correct Ans"
$ws.Cells.Item(15,4).Value = "Student Ans"
$ws.Cells.Item(15,5).Value = "Correct Ans"
Format-Bold $ws.Cells.Item(15,1)
Format-Bold $ws.Cells.Item(15,2)
Format-Bold $ws.Cells.Item(15,4)
Format-Bold $ws.Cells.Item(15,5)

# ---- left answer table : rows 16-40, columns A (student) / B (correct) ---------
$leftAnswers = @(
    @("Option D","Option A"),
    @("Option D","Option D"),
    @("Option B","Option B"),
    @("Option C","Option C"),
    @("Option B","Option B"),
    @("Option B","Option C"),
    @("Option A","Option D"),
    @("","Option D"),
    @("","Option A"),
    @("Option A","Option A"),
    @("Option C","Option C"),
    @("Option A","Option A"),
    @("Option D","Option D"),
    @("Option D","Option D"),
    @("Option B","Option B"),
    @("Option C","Option D"),
    @("Option C","Option C"),
    @("Option D","Option D"),
    @("Option B","Option B"),
    @("","Option D"),
    @("Option A","Option A"),
    @("Option A","Option A"),
    @("Option B","Option A"),
    @("Option D","Option D"),
    @("Option A","Option D")
)

$row = 16
foreach ($pair in $leftAnswers) {
    $studentAns = $pair[0]
    $correctAns = $pair[1]

    $studentCell = $ws.Cells.Item($row,1)
    $correctCell = $ws.Cells.Item($row,2)

    $studentCell.Value = $studentAns
    $correctCell.Value = $correctAns

    if ($studentAns -eq $correctAns) {
        Format-Colored $studentCell $GREEN
    } else {
        Format-Colored $studentCell $RED
    }
    Format-Colored $correctCell $BLUE

    $row++
}

# ---- right answer table : rows 16-18, columns D (student) / E (correct) --------
$rightAnswers = @(
    @("Option C","Option A"),
    @("Option C","Option C"),
    @("Option D","Option D")
)

$row = 16
foreach ($pair in $rightAnswers) {
    $studentAns = $pair[0]
    $correctAns = $pair[1]

    $studentCell = $ws.Cells.Item($row,4)
    $correctCell = $ws.Cells.Item($row,5)

    $studentCell.Value = $studentAns
    $correctCell.Value = $correctAns

    if ($studentAns -eq $correctAns) {
        Format-Colored $studentCell $GREEN
    } else {
        Format-Colored $studentCell $RED
    }
    Format-Colored $correctCell $BLUE

    $row++
}
